$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2063233333333333
$ws.Range("H2").Value = 0.61897
$ws.Range("M2").Value = 14.440165
$ws.Range("N2").Value = 43.320495
$ws.Range("O2").Value = 0.1441015470002482
$ws.Range("P2").Value = 0.1441015470002482
$ws.Range("Q2").Value = 2.979342976683333
$ws.Range("R2").Value = 26.81408679015
$ws.Range("S2").Value = 0.1441015470002482
$ws.Range("T2").Value = 0.1441015470002482

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2063233333333333
$ws.Range("H3").Value = 0.61897
$ws.Range("O3").Value = 0.3846359116098663
$ws.Range("P3").Value = 0.3846359116098662
$ws.Range("Q3").Value = 7.952463562608889
$ws.Range("R3").Value = 71.57217206348
$ws.Range("S3").Value = 0.3846359116098663
$ws.Range("T3").Value = 0.3846359116098662

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2063233333333333
$ws.Range("H4").Value = 0.61897
$ws.Range("M4").Value = 21.954262
$ws.Range("N4").Value = 65.862786
$ws.Range("O4").Value = 0.2190863551385157
$ws.Range("P4").Value = 0.2190863551385156
$ws.Range("Q4").Value = 4.529676516713333
$ws.Range("R4").Value = 40.76708865042
$ws.Range("S4").Value = 0.2190863551385157
$ws.Range("T4").Value = 0.2190863551385156

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2063233333333333
$ws.Range("H5").Value = 0.61897
$ws.Range("M5").Value = 25.27013633333333
$ws.Range("N5").Value = 75.81040899999999
$ws.Range("O5").Value = 0.2521761862513699
$ws.Range("P5").Value = 0.2521761862513699
$ws.Range("Q5").Value = 5.21381876208111
$ws.Range("R5").Value = 46.92436885873
$ws.Range("S5").Value = 0.2521761862513699
$ws.Range("T5").Value = 0.2521761862513699

$wb.Save()
